# The table in this workbook lists EPP status-code test cases. Historically it
# included, for every status code, a redundant pair of "add" rows (one
# expecting success, one expecting the server to reject adding a status the
# object already has) and for "rem" tests a redundant pair (one expecting
# success, one expecting the server to reject removing a status the object
# doesn't have). Per the commit message, those redundant "reject-because-
# already-present / reject-because-not-present" checks should be removed,
# except for the ones that were already exercising a *different* status
# (clientHold for "add", invalidStatusCode for "rem").
#
# The simplest, most reliable way to reproduce the resulting table is to
# write the final data directly into B10:E32 (the surviving rows, in their
# final order) and then delete the now-unused trailing rows (33:41), which
# also shrinks the Table/AutoFilter/dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("add", "clientDeleteProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("add", "clientHold", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("add", "clientRenewProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("add", "clientTransferProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("add", "clientUpdateProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("add", "clientHold", "fail", "EPP_UNEXPECTED_COMMAND_SUCCESS"),
    @("rem", "clientUpdateProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("add", "linked", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "ok", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "pendingCreate", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "pendingDelete", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "pendingTransfer", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "pendingUpdate", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "serverDeleteProhibited", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "serverRenewProhibited", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "serverTransferProhibited", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "serverUpdateProhibited", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("add", "invalidStatusCode", "fail", "EPP_DOMAIN_UPDATE_SERVER_ACCEPTS_INVALID_STATUS_CODE"),
    @("rem", "clientDeleteProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("rem", "clientHold", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("rem", "clientRenewProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("rem", "clientTransferProhibited", "pass", "EPP_UNEXPECTED_COMMAND_FAILURE"),
    @("rem", "invalidStatusCode", "fail", "EPP_UNEXPECTED_COMMAND_SUCCESS")
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("B$r").Value = $data[0]
    $ws.Range("C$r").Value = $data[1]
    $ws.Range("D$r").Value = $data[2]
    $ws.Range("E$r").Value = $data[3]
}

$lastNewRow = $startRow + $rows.Count - 1

# Remove the now-superfluous trailing rows (this also shrinks the
# Table/AutoFilter ref and the sheet dimension automatically).
$oldLastRow = 41
$ws.Range("A" + ($lastNewRow + 1) + ":A" + $oldLastRow).EntireRow.Delete()

$ws.Range("C" + $lastNewRow).Select()
